# Acknowledge page fix and plot report background
#
# The crop-yield "outlier" acknowledgement screen was incorrectly reusing
# the `ack_ok` session variable (and its branch-fail/clause formulas) that
# belongs to the crop-height outlier screen. Give the yield screen its own
# `ack_ok_2` variable (and matching message / formulas), and register the
# new session variable on the "model" sheet.

$wb = $excel.ActiveWorkbook
$wsSurvey = $wb.Worksheets.Item("survey")
$wsModel = $wb.Worksheets.Item("model")

# ---------------------------------------------------------------------
# 1. "survey" sheet: point the yield-outlier acknowledge block at the new
#    ack_ok_2 variable instead of the (height-screen) ack_ok variable.
# ---------------------------------------------------------------------

# acknowledge prompt's display text (yield-specific wording)
$wsSurvey.Range("I59").Value2 = "The yield you have entered is an outlier.  Are you certain of this measurement?"

# acknowledge prompt's session-variable name
$wsSurvey.Range("E59").Value2 = "ack_ok_2"

# acknowledge prompt's clause/condition formula
$wsSurvey.Range("K59").Value2 = "data('ack_ok_2') || calculates.below_max_yield()"

# branch_fail condition formula, right below
$wsSurvey.Range("C60").Value2 = "! data('ack_ok_2')"

# ---------------------------------------------------------------------
# 2. "model" sheet: register the new ack_ok_2 session variable.
#    Insert a row right after the existing ack_ok row so the table reads
#    ack_ok, ack_ok_2, max_height_disp, max_yield_disp, ...
# ---------------------------------------------------------------------
$wsModel.Rows("3:3").Insert()
$wsModel.Range("A3").Value2 = "ack_ok_2"
$wsModel.Range("B3").Value2 = "boolean"
$wsModel.Range("C3").Value2 = $true

# ---------------------------------------------------------------------
# 3. Cosmetic view-state updates captured in the workbook (selection /
#    active cell on the "model" and "survey" sheets). Select the "model"
#    sheet's cell first so the final active tab/selection ends back on
#    "survey", matching the saved workbook state.
# ---------------------------------------------------------------------
$wsModel.Range("A43").Select()
$wsSurvey.Range("E17").Select()
